$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (codevita_exam_solutions) updates
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = "2026-02-18T20:03:57.516403+00:00"
$ws.Range("H3").Value = 24
$ws.Range("L3").Value = "[486988, 486982, 487002, 487051, 487036, 487065, 487059, 487110, 487103, 487102, 487114, 487074, 487066, 487055, 487138, 487137, 487149, 487164, 487144, 487159, 487228, 487218, 487227, 487257]"
